$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.412.15"
$ws.Range("E2").Value = "  -2.17%  "

$ws.Range("D3").Value = "1.654.24"

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "213.72"
$ws.Range("E5").Value = "  -1.36%  "

$ws.Range("D6").Value = "0.511"
$ws.Range("E6").Value = "  -1.41%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").Value = "24.20"
$ws.Range("E8").Value = "  +0.50%  "

$ws.Range("D9").Value = "0.261"
$ws.Range("E9").Value = "  -1.16%  "

$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("E11").Value = "  -0.72%  "

$ws.Range("D12").Value = "1.887.86"
$ws.Range("E12").Value = "  -1.93%  "

$ws.Range("D13").Value = "1.647.47"
$ws.Range("E13").Value = "  -2.35%  "

$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").Value = "4.09"
$ws.Range("E14").Value = "  -2.27%  "

$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "0.572"
$ws.Range("E15").Value = "  +2.32%  "

$ws.Range("D16").Value = "65.91"
$ws.Range("E16").Value = "  -1.40%  "

$ws.Range("D17").Value = "27.412.21"
$ws.Range("E17").Value = "  -1.93%  "

$ws.Range("D18").Value = "234.12"
$ws.Range("E18").Value = "  -6.41%  "

$ws.Range("E19").Value = "  -1.90%  "

$ws.Range("E20").Value = "  -2.77%  "

$ws.Range("D21").Value = "0.999"
$ws.Range("E21").Value = "  -0.13%  "

$ws.Range("D22").Value = "4.39"
$ws.Range("E22").Value = "  -2.88%  "

$ws.Range("E23").Value = "  -2.34%  "

$ws.Range("D25").Value = "147.07"
$ws.Range("E25").Value = "  -0.17%  "

$ws.Range("D26").Value = "7.21"
$ws.Range("E26").Value = "  -1.73%  "

$ws.Range("D27").Value = "16.00"
$ws.Range("E27").Value = "  -2.93%  "

$ws.Range("E28").Value = "  -0.21%  "

$ws.Range("E29").Value = "  -1.96%  "

$ws.Range("D30").Value = "0.0497"
$ws.Range("E30").Value = "  -1.11%  "

$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  -5.04%  "

$ws.Range("E32").Value = "  -2.28%  "

$ws.Range("D33").Value = "1.462.90"
$ws.Range("E33").Value = "  +2.61%  "

$ws.Range("E34").Value = "  -2.28%  "

$ws.Range("E36").Value = "  -0.76%  "

$ws.Range("D37").Value = "0.911"
$ws.Range("E37").Value = "  -3.43%  "

$ws.Range("E38").Value = "  -3.35%  "

$ws.Range("E39").Value = "  -1.31%  "

$ws.Range("E40").Value = "  -0.31%  "

$ws.Range("E41").Value = "  -0.14%  "

$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("E43").Value = "  -5.47%  "

$ws.Range("E44").Value = "  -0.46%  "

$ws.Range("B45").Value = "RocketPoolETH"
$ws.Range("C45").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D45").Value = "1.796.22"
$ws.Range("E45").Value = "  -1.97%  "

$ws.Range("B46").Value = "TrustWalletToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D46").Value = "0.784"
$ws.Range("E46").Value = "  -1.54%  "

$ws.Range("E47").Value = "  +0.28%  "

$ws.Range("D48").Value = "88.35"
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("D49").Value = "0.0₆0105"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("E50").Value = "  -1.49%  "

$ws.Range("D51").Value = "7.80"
$ws.Range("E51").Value = "  -1.41%  "

